# May 9th changes: a new sample (x,y,z) was recorded at the start of the
# series (becoming row 2, pushing the previous rows 2-21 down to 3-22),
# and 9 brand-new samples were appended at the end (rows 23-31).
#
# Rather than using EntireRow.Insert() (which, in this runtime, stamps the
# freshly inserted row with a copy of the header row's bold style - leaving
# an extraneous, unused cell-format entry behind even after ClearFormats),
# we just rewrite the whole A2:C31 block directly with its final values.
# The net effect on the sheet is identical to "insert one row up top, then
# append nine more at the bottom".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(-0.03398093824483917,  -2.852332382786031,    0.3628882449500415)
    3  = @(-0.1266735037978812,   -1.426815969603404,    0.2018442184341198)
    4  = @(-0.1538507725511279,    0.2260856117521008,  -0.09954921262604832)
    5  = @(0.4955612986671659,     3.062298653077098,   -0.5691612070920504)
    6  = @(1.454008883359481,      4.839717144868811,   -0.2820296159812385)
    7  = @(1.137240985218359,      3.374921506764936,    0.5055315214760452)
    8  = @(0.2921525373750798,     1.951630319867816,    0.7453397719227535)
    9  = @(0.1624465198541174,     1.506097605033798,    0.4439899596024537)
    10 = @(0.2823693378239257,    -0.7324928641319217,  -0.04955176400894934)
    11 = @(-0.6597999164036341,   -4.261745044163296,   -0.7145596061434063)
    12 = @(-1.235334702292262,    -4.321111241165461,   -0.293212206996212)
    13 = @(0.393702644170544,     -1.912042505887086,    0.653996666171116)
    14 = @(0.5366638071683012,    -2.253858975001746,    0.0269684557403839)
    15 = @(-0.129132547548839,    -0.7719840942596894,   0.03628414990950613)
    16 = @(0.3077981770038601,     1.624053824921043,    0.1965552446793536)
    17 = @(1.183255352536026,      3.720431172117896,    0.4516974523359422)
    18 = @(0.7590655258723678,     3.928779942648755,    1.171988606452941)
    19 = @(0.2291679642334281,     2.087977978647975,    1.354070066189282)
    20 = @(0.03184602683296017,    0.9907392433711465,   0.4831009315592913)
    21 = @(-0.04444044737183325,   0.1521864691559109,  -0.2139835976520377)
    22 = @(-0.4428928944529248,   -2.088161782342548,   -1.22169929499528)
    23 = @(-1.581159264457474,    -3.781517471585969,   -2.15734222470499)
    24 = @(-0.5357818153439736,   -0.6552340047700065,   0.5789350879435637)
    25 = @(-0.4312272305999488,   -0.03695735122476339, -0.2086323031357354)
    26 = @(0.1521366113910867,     0.3846518628451288,  -0.3579327458021597)
    27 = @(0.3417635331956719,     0.5353018106246487,   0.2242374224018085)
    28 = @(-0.1421539567563001,    0.2858568746216445,  -0.08254160519157175)
    29 = @(-0.06803667803808118,   0.1790894811250734,   0.08434615633925557)
    30 = @(-0.08848196070413178,  -0.1032362286837731,   0.232543302129726)
    31 = @(0.06249837318853448,    0.04699299066346516,  0.2259266389419836)
}

for ($r = 2; $r -le 31; $r++) {
    $vals = $data[$r]
    $ws.Range("A$r").Value = $vals[0]
    $ws.Range("B$r").Value = $vals[1]
    $ws.Range("C$r").Value = $vals[2]
}
